# Edit script: "Endringer i figurer mm. Omgjoere til SVG-format"
# Inserts a new sheet "inntekt_mottakere" between "mottaker" and "husholdning",
# populates it with data, and updates selections/active tab accordingly.

$wb = $excel.ActiveWorkbook

$wsMottaker = $wb.Worksheets.Item("mottaker")
$wsHusholdning = $wb.Worksheets.Item("husholdning")

# --- 1. Update the "mottaker" sheet's selection (no longer the active tab) ---
$wsMottaker.Activate()
$wsMottaker.Range("C1").Select()

# --- 2. Insert the new sheet right after "mottaker" ---
$wsNew = $wb.Worksheets.Add($null, $wsMottaker)
$wsNew.Name = "inntekt_mottakere"

# --- 3. Header row (copy formatting + shared strings from the "mottaker" header) ---
$wsMottaker.Range("A1:D1").Copy($wsNew.Range("A1"))

# --- 4. Data rows 2-29 : gr / ar / inntekt / kg ---
$data = @(
  @("[0.85,0.95)", 2015, 0.22618185889788001, 0.85, 0),
  @("[0.85,0.95)", 2016, 0.25317161389929438, 0.85, 0),
  @("[0.85,0.95)", 2017, 0.27495898841839328, 0.85, 0),
  @("[0.85,0.95)", 2018, 0.280119165536863, 0.85, 0),
  @("[0.85,0.95)", 2019, 0.2882241062491474, 0.85, 0),
  @("[0.85,0.95)", 2020, 0.27701246795345991, 0.85, 0),
  @("[0.85,0.95)", 2021, 0.27638685550139391, 0.85, 0),
  @("[0.90,0.95)", 2015, 0.182315, 0.9, 1),
  @("[0.90,0.95)", 2016, 0.211671, 0.9, 1),
  @("[0.90,0.95)", 2017, 0.23494100000000001, 0.9, 1),
  @("[0.90,0.95)", 2018, 0.23594100000000001, 0.9, 1),
  @("[0.90,0.95)", 2019, 0.23769999999999999, 0.9, 1),
  @("[0.90,0.95)", 2020, 0.228655, 0.9, 1),
  @("[0.90,0.95)", 2021, 0.23419100000000001, 0.9, 1),
  @("[0.95,~)", 2015, 0.14872290341299491, 0.9, 1),
  @("[0.95,~)", 2016, 0.19353657244446279, 0.9, 1),
  @("[0.95,~)", 2017, 0.21843568620692241, 0.9, 1),
  @("[0.95,~)", 2018, 0.2292012626661232, 0.9, 1),
  @("[0.95,~)", 2019, 0.23094830701794189, 0.9, 1),
  @("[0.95,~)", 2020, 0.22084969934054799, 0.9, 1),
  @("[0.95,~)", 2021, 0.22815654902582411, 0.9, 1),
  @("[0.95,~)", 2015, 0.14872290341299491, 0.85, 0),
  @("[0.95,~)", 2016, 0.19353657244446279, 0.85, 0),
  @("[0.95,~)", 2017, 0.21843568620692241, 0.85, 0),
  @("[0.95,~)", 2018, 0.2292012626661232, 0.85, 0),
  @("[0.95,~)", 2019, 0.23094830701794189, 0.85, 0),
  @("[0.95,~)", 2020, 0.22084969934054799, 0.85, 0),
  @("[0.95,~)", 2021, 0.22815654902582411, 0.85, 0)
)

$r = 2
foreach ($row in $data) {
    $wsNew.Cells.Item($r, 1).Value = $row[0]
    $wsNew.Cells.Item($r, 2).Value = $row[1]
    $cC = $wsNew.Cells.Item($r, 3)
    $cC.Value = $row[2]
    $cC.NumberFormat = "#,##0.00"
    $cD = $wsNew.Cells.Item($r, 4)
    $cD.Value = $row[3]
    if ($row[4] -eq 1) {
        $cD.NumberFormat = "#,##0.00"
    }
    $r++
}

# --- 5. Leftover formatting-only cells (columns E-H) carried over from the source data ---
$extraFormats = @(
  @("F2", "general"),
  @("G2", "general"),
  @("H2", "general"),
  @("E3", "general"),
  @("F3", "general"),
  @("G3", "general"),
  @("H3", "general"),
  @("E4", "general"),
  @("F4", "general"),
  @("G4", "general"),
  @("H4", "general"),
  @("E5", "general"),
  @("F5", "general"),
  @("G5", "general"),
  @("H5", "general"),
  @("E6", "general"),
  @("F6", "general"),
  @("G6", "general"),
  @("H6", "general"),
  @("E7", "general"),
  @("F7", "general"),
  @("G7", "general"),
  @("H7", "general"),
  @("E8", "general"),
  @("F8", "general"),
  @("G8", "general"),
  @("H8", "general"),
  @("E9", "general"),
  @("F9", "#,##0.0000"),
  @("G9", "general"),
  @("H9", "general"),
  @("E10", "general"),
  @("F10", "#,##0.0000"),
  @("G10", "general"),
  @("H10", "general"),
  @("E11", "general"),
  @("F11", "#,##0.0000"),
  @("G11", "general"),
  @("H11", "general"),
  @("E12", "general"),
  @("F12", "#,##0.0000"),
  @("G12", "general"),
  @("H12", "general"),
  @("E13", "general"),
  @("F13", "#,##0.0000"),
  @("G13", "general"),
  @("H13", "general"),
  @("E14", "general"),
  @("F14", "#,##0.0000"),
  @("G14", "general"),
  @("H14", "general"),
  @("E15", "general"),
  @("F15", "#,##0.0000"),
  @("G15", "general"),
  @("H15", "general"),
  @("E16", "general"),
  @("F16", "#,##0.00"),
  @("E17", "general"),
  @("F17", "#,##0.00"),
  @("G17", "general"),
  @("E18", "general"),
  @("F18", "#,##0.00"),
  @("G18", "general"),
  @("E19", "general"),
  @("F19", "#,##0.00"),
  @("G19", "general"),
  @("E20", "general"),
  @("F20", "#,##0.00"),
  @("E21", "#,##0.00"),
  @("F21", "#,##0.00"),
  @("E22", "#,##0.00"),
  @("F22", "#,##0.00")
)

foreach ($item in $extraFormats) {
    $wsNew.Range($item[0]).NumberFormat = $item[1]
}

# --- 6. Column width for column A ---
$wsNew.Columns.Item(1).ColumnWidth = 10.42

# --- 7. Sheet view: zoom + selection, make this the active tab ---
$wsNew.Activate()
$excel.ActiveWindow.Zoom = 130
$wsNew.Range("C13").Select()

Write-Output "Edit complete"
